$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the two missing cells for the existing last row (row 10)
$ws.Range("X10").Value = -3.4100040000000149
$ws.Range("Y10").Value = "Down"

# Append the new traded row (row 11)
$ws.Range("A11").Value = 42654.886666666665
$ws.Range("B11").Value = 12
$ws.Range("C11").Value = "Buy"
$ws.Range("D11").Value = 30
$ws.Range("E11").Value = 22664
$ws.Range("F11").Value = 3961
$ws.Range("G11").Value = 60
$ws.Range("H11").Value = 35
$ws.Range("I11").Value = 88
$ws.Range("J11").Value = 11
$ws.Range("K11").Value = 38421
$ws.Range("L11").Value = 338
$ws.Range("M11").Value = 198
$ws.Range("N11").Value = 84
$ws.Range("O11").Value = 11
$ws.Range("P11").Value = "Noun"
$ws.Range("Q11").Value = 39.313912976930268
$ws.Range("R11").Value = 1.8
$ws.Range("S11").Value = 0.0864
$ws.Range("T11").Value = -0.0115
$ws.Range("U11").Value = 5.85
$ws.Range("V11").Value = "N/A"
$ws.Range("W11").Value = 0

# Match the existing number formats (date style from A column, percent style from S/T columns)
# by copying formatting from the row above, instead of minting new custom number formats.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)

$ws.Range("S10").Copy()
$ws.Range("S11:T11").PasteSpecial(-4122)
